$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the hourly crypto price/volume snapshot (GitHub Actions scheduled
# update). Only the numeric Price (column D) and Volume 1h % (column E) cells
# for the affected coins change, except for rows 14-15 where Chainlink and
# Polygon swap ranking positions (their Coin name, Link, Price and Volume all
# move together).
#
# Price/volume values are plain-text cells in the source sheet (e.g. prices use
# "." as a thousands separator, like "35.454.48", and volumes keep padding
# spaces, like "  +2.83%  "). Some of the new values would otherwise look like
# an ordinary number to Excel (e.g. "231.30") and get silently reinterpreted as
# a numeric value (dropping the trailing zero, switching to scientific notation,
# etc.), so for those cells we briefly force a Text number format while writing
# the value, then restore the original (default/"Normal") cell style so no
# visible formatting changes are left behind.

# Row 2: Bitcoin
$ws.Range("D2").Value = "35.454.48"
$ws.Range("E2").Value = "  +2.83%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.838.45"
$ws.Range("E3").Value = "  +1.80%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.30%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.91%  "

# Row 6: XRP
$ws.Range("E6").Value = "  +1.49%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.35%  "

# Row 8: Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.28%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  +7.26%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0704"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.25%  "

# Row 11: TRON
$ws.Range("E11").Value = "  +2.50%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.104.81"
$ws.Range("E12").Value = "  +1.81%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.840.03"
$ws.Range("E13").Value = "  +1.90%  "

# Row 14: Chainlink
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.673"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.66%  "

# Row 15: Polygon
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "11.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.58%  "

# Row 16: Polkadot
$ws.Range("E16").Value = "  +7.57%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "35.393.28"
$ws.Range("E17").Value = "  +2.73%  "

# Row 18: Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.88%  "

# Row 19: ShibaInu
$ws.Range("D19").Value = "0.0₃0801"
$ws.Range("E19").Value = "  +4.28%  "

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.25%  "

# Row 21: Avalanche
$ws.Range("E21").Value = "  +8.23%  "

# Row 22: Uniswap
$ws.Range("E22").Value = "  +14.70%  "

# Row 23: Dai
$ws.Range("E23").Value = "  +0.31%  "

# Row 24: Toncoin
$ws.Range("E24").Value = "  +0.90%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.64"
$ws.Range("D25").Style = "Normal"

# Row 26: Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.13%  "

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "

# Row 28: Stellar
$ws.Range("E28").Value = "  -0.78%  "

# Row 29: PancakeSwap
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +27.58%  "

# Row 30: BinanceUSD
$ws.Range("E30").Value = "  +0.28%  "

# Row 31: EURNeutrino
$ws.Range("D31").Value = "3.356.50"
$ws.Range("E31").Value = "  +38.15%  "

# Row 32: Hedera
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0553"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.50%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +6.28%  "

# Row 34: Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.54%  "

# Row 35: LidoDAOToken
$ws.Range("E35").Value = "  +1.56%  "

# Row 36: Aave
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "95.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.52%  "

# Row 37: ImmutableX
$ws.Range("E37").Value = "  +7.21%  "

# Row 38: TrustWalletToken
$ws.Range("E38").Value = "  +5.86%  "

# Row 39: InjectiveProtocol
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.84%  "

# Row 40: Maker
$ws.Range("D40").Value = "1.347.86"
$ws.Range("E40").Value = "  +3.12%  "

# Row 41: RenderToken
$ws.Range("E41").Value = "  +5.09%  "

# Row 42: VeChain
$ws.Range("E42").Value = "  +4.90%  "

# Row 43: ARBITRUM
$ws.Range("E43").Value = "  +5.78%  "

# Row 44: WEMIXToken
$ws.Range("E44").Value = "  +4.58%  "

# Row 45: HuobiToken
$ws.Range("E45").Value = "  +0.87%  "

# Row 46: MXToken
$ws.Range("E46").Value = "  +0.63%  "

# Row 47: FraxShare
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.63%  "

# Row 48: Kaspa
$ws.Range("E48").Value = "  +1.42%  "

# Row 49: RocketPoolETH
$ws.Range("D49").Value = "2.007.71"
$ws.Range("E49").Value = "  +2.03%  "

# Row 50: PaxDollar
$ws.Range("E50").Value = "  +0.43%  "

# Row 51: Quant
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.68%  "
